# Apply cell updates from the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force Excel to store the value as literal text instead of
    # auto-converting number-looking strings (e.g. "498.09") to a
    # floating point number, then restore the default (unstyled) cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "54.097.54"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.268.40"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "498.09"
$ws.Range("E5").Value = "  -0.09%  "
Set-TextValue $ws.Range("D6") "127.81"
$ws.Range("E6").Value = "  -0.45%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue $ws.Range("D8") "0.525"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +0.38%  "
Set-TextValue $ws.Range("D11") "0.335"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "2.670.70"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "54.054.03"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.286.85"
$ws.Range("E17").Value = "  -0.75%  "
Set-TextValue $ws.Range("D18") "10.21"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  +1.40%  "
Set-TextValue $ws.Range("D20") "302.30"
$ws.Range("E20").Value = "  -1.32%  "
Set-TextValue $ws.Range("D21") "6.31"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  +0.28%  "
Set-TextValue $ws.Range("D23") "60.99"
$ws.Range("E23").Value = "  -3.36%  "
Set-TextValue $ws.Range("D24") "1.00"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -1.81%  "
Set-TextValue $ws.Range("D26") "7.27"
$ws.Range("E26").Value = "  +1.43%  "
Set-TextValue $ws.Range("D27") "172.91"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D28") "1.60"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0689"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("E34").Value = "  +0.00%  "
Set-TextValue $ws.Range("D35") "0.934"
$ws.Range("E35").Value = "  +7.68%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  +0.56%  "
Set-TextValue $ws.Range("D38") "0.371"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D41") "124.96"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "4.79"
$ws.Range("E42").Value = "  -1.42%  "
Set-TextValue $ws.Range("D43") "0.0490"
$ws.Range("E43").Value = "  +1.15%  "
Set-TextValue $ws.Range("D44") "0.0888"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  -1.15%  "
Set-TextValue $ws.Range("D46") "238.01"
$ws.Range("E46").Value = "  -2.72%  "
Set-TextValue $ws.Range("D47") "0.370"
$ws.Range("E47").Value = "  -1.34%  "
Set-TextValue $ws.Range("D48") "0.0204"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.37%  "
Set-TextValue $ws.Range("D50") "16.13"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  -0.53%  "
